$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5069
    $ws.Range("F3").Value = 155
    $ws.Range("F4").Value = 887
}
